# Apply updated crypto price/volume data pulled from the latest feed,
# including the two pairs of rows that swapped rank order (21/22 and 26/27).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.187.37'
$ws.Range('E2').Value = '  -2.91%  '
$ws.Range('D3').Value = '1.850.73'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'0.7035"
$ws.Range('E5').Value = '  -4.69%  '
$ws.Range('D6').Value = "'239.11"
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = "'0.3053"
$ws.Range('E8').Value = '  -3.73%  '
$ws.Range('D9').Value = "'0.07431"
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('D10').Value = "'23.40"
$ws.Range('E10').Value = '  -5.65%  '
$ws.Range('D11').Value = "'0.08158"
$ws.Range('D12').Value = '1.893.60'
$ws.Range('E12').Value = '  +1.42%  '
$ws.Range('D13').Value = "'0.7284"
$ws.Range('E13').Value = '  -3.88%  '
$ws.Range('D14').Value = "'5.218"
$ws.Range('E14').Value = '  -3.32%  '
$ws.Range('D15').Value = "'89.20"
$ws.Range('D16').Value = '29.431.39'
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').Value = "'5.786"
$ws.Range('E17').Value = '  -5.92%  '
$ws.Range('D18').Value = "'238.96"
$ws.Range('E18').Value = '  -4.42%  '
$ws.Range('D19').Value = "'13.11"
$ws.Range('E19').Value = '  -3.30%  '
$ws.Range('D20').Value = "'0.000007672"
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.147.27'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = "'1.001"
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = "'7.602"
$ws.Range('E24').Value = '  -3.69%  '
$ws.Range('D25').Value = "'9.018"
$ws.Range('E25').Value = '  -2.90%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'160.76"
$ws.Range('E26').Value = '  -1.59%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.1456"
$ws.Range('E27').Value = '  -6.64%  '
$ws.Range('D28').Value = "'18.12"
$ws.Range('E28').Value = '  -3.16%  '
$ws.Range('D29').Value = "'1.979"
$ws.Range('E29').Value = '  -3.46%  '
$ws.Range('D30').Value = "'1.411"
$ws.Range('E30').Value = '  -3.77%  '
$ws.Range('D31').Value = "'4.524"
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = "'1.493"
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('E33').Value = '  -4.69%  '
$ws.Range('D34').Value = "'0.05206"
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('D35').Value = "'1.188"
$ws.Range('E35').Value = '  -5.06%  '
$ws.Range('D36').Value = "'1.037"
$ws.Range('E36').Value = '  +3.68%  '
$ws.Range('D37').Value = "'0.7082"
$ws.Range('E37').Value = '  -7.84%  '
$ws.Range('D38').Value = "'2.661"
$ws.Range('E38').Value = '  -2.18%  '
$ws.Range('D39').Value = "'0.01869"
$ws.Range('E39').Value = '  -4.75%  '
$ws.Range('D41').Value = "'0.9420"
$ws.Range('E41').Value = '  +8.04%  '
$ws.Range('D42').Value = "'6.046"
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = "'0.4311"
$ws.Range('E43').Value = '  -5.86%  '
$ws.Range('D44').Value = '1.068.24'
$ws.Range('E44').Value = '  -2.10%  '
$ws.Range('D45').Value = "'70.60"
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = "'103.67"
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').Value = '2.029.78'
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('E50').Value = '  -7.06%  '
$ws.Range('D51').Value = "'9.159"
$ws.Range('E51').Value = '  -4.08%  '
